$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("E3").Value = "  -4.58%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("E6").Value = "  -7.06%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -10.66%  "
$ws.Range("E9").Value = "  -4.47%  "
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("E14").Value = "  -6.21%  "
$ws.Range("E15").Value = "  -4.74%  "
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("E17").Value = "  -4.51%  "
$ws.Range("E18").Value = "  -4.74%  "
$ws.Range("E19").Value = "  -3.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "316.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("E22").Value = "  -6.15%  "
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  +10.73%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -4.67%  "
$ws.Range("E29").Value = "  -9.47%  "
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("E31").Value = "  -5.74%  "
$ws.Range("E32").Value = "  -7.26%  "
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("E34").Value = "  -4.78%  "
$ws.Range("E35").Value = "  -4.74%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -4.61%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("E40").Value = "  -9.06%  "
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "140.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("E46").Value = "  -8.99%  "
$ws.Range("E48").Value = "  -4.58%  "
$ws.Range("E49").Value = "  -11.58%  "
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("E51").Value = "  -3.62%  "
